# [#Tutorial 17] Add row data table from UI
# Adds a new "Sheet2" worksheet (placed right after Sheet1) containing a
# small Name/Age table, then makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet immediately after Sheet1 so the tab order is
# Sheet1, Sheet2 (Worksheets.Add() alone would insert before the active
# sheet, i.e. at the front).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Age"

# Data rows
$ws2.Range("A2").Value = "Mane"
$ws2.Range("B2").Value = 27

$ws2.Range("A3").Value = "Virgil"
$ws2.Range("B3").Value = 30

$ws2.Range("A4").Value = "Firmino"
$ws2.Range("B4").Value = 22

$ws2.Range("A5").Value = "Salah"
$ws2.Range("B5").Value = 45

# Make the new sheet the active tab / selection, matching the UI action
# that created it.
$ws2.Select()
$ws2.Range("G7").Select()
